$wb = $excel.ActiveWorkbook


# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 5
$ws.Range("H5").Value = 56.8
$ws.Range("I5").Value = 46.25
$ws.Range("K5").Value = 46.25
$ws.Range("M5").Value = 68.75
# Row 40
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
# Row 41
$ws.Range("H41").Value = 2208.0557
$ws.Range("I41").Value = 1596.875
$ws.Range("K41").Value = 1596.875
$ws.Range("M41").Value = -1156.875
# Row 64
$ws.Range("H64").Value = 4236.091
$ws.Range("I64").Value = 3437.125
$ws.Range("K64").Value = 3437.125
$ws.Range("M64").Value = -3189.125
# Row 67
$ws.Range("H67").Value = 4236.091
$ws.Range("I67").Value = 3437.125
$ws.Range("K67").Value = 3437.125
$ws.Range("M67").Value = -2579.125
# Row 74
$ws.Range("H74").Value = 5229
$ws.Range("I74").Value = 4732.6665
$ws.Range("K74").Value = 4732.6665
$ws.Range("M74").Value = -3796.6665
# Row 77
$ws.Range("H77").Value = 5229
$ws.Range("I77").Value = 4732.6665
$ws.Range("K77").Value = 23663.3325
$ws.Range("M77").Value = -18983.3325
# Row 98
$ws.Range("H98").Value = 992.7308
$ws.Range("I98").Value = 1078.5
$ws.Range("J98").Value = 706.8333
$ws.Range("K98").Value = 1078.5
$ws.Range("L98").Value = 706.8333
$ws.Range("M98").Value = 419.5
$ws.Range("N98").Value = -3702.8333
# Row 113
$ws.Range("H113").Value = 5400.3335
$ws.Range("I113").Value = 4241.6
$ws.Range("J113").Value = 6848.75
$ws.Range("K113").Value = 4241.6
$ws.Range("L113").Value = 6848.75
$ws.Range("M113").Value = -987.6000000000004
$ws.Range("N113").Value = -13356.75
# Row 122
$ws.Range("H122").Value = 992.7308
$ws.Range("I122").Value = 1078.5
$ws.Range("J122").Value = 706.8333
$ws.Range("K122").Value = 3235.5
$ws.Range("L122").Value = 2120.4999
$ws.Range("M122").Value = -785.5
$ws.Range("N122").Value = -7020.4999
# Row 132
$ws.Range("H132").Value = 1611.6111
$ws.Range("I132").Value = 1611.6111
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4834.8333
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2304.8333
$ws.Range("N132").ClearContents()
# Row 133
$ws.Range("H133").Value = 95646.664
$ws.Range("J133").Value = 95646.664
$ws.Range("L133").Value = 95646.664
$ws.Range("N133").Value = -105766.664
# Row 139
$ws.Range("H139").Value = 99990
$ws.Range("J139").Value = 99990
$ws.Range("L139").Value = 99990
$ws.Range("N139").Value = -110270
# Row 140
$ws.Range("H140").Value = 91990
$ws.Range("J140").Value = 91990
$ws.Range("L140").Value = 91990
$ws.Range("N140").Value = -102350

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 21
$ws.Range("H21").Value = 799.375
$ws.Range("I21").Value = 798.3333
$ws.Range("K21").Value = 798.3333
$ws.Range("M21").Value = -424.3333
# Row 32
$ws.Range("H32").Value = 8625.459999999999
$ws.Range("I32").Value = 6363.6045
$ws.Range("J32").Value = 22519.715
$ws.Range("K32").Value = 6363.6045
$ws.Range("L32").Value = 22519.715
$ws.Range("M32").Value = -6076.6045
$ws.Range("N32").Value = -23093.715
# Row 61
$ws.Range("H61").Value = 3065.3823
$ws.Range("I61").Value = 2912.0667
$ws.Range("J61").Value = 4215.25
$ws.Range("K61").Value = 2912.0667
$ws.Range("L61").Value = 4215.25
$ws.Range("M61").Value = -2700.0667
$ws.Range("N61").Value = -4639.25
# Row 74
$ws.Range("H74").Value = 2048.6667
$ws.Range("I74").Value = 1598.2858
$ws.Range("K74").Value = 1598.2858
$ws.Range("M74").Value = -724.2858000000001
# Row 77
$ws.Range("H77").Value = 2048.6667
$ws.Range("I77").Value = 1598.2858
$ws.Range("K77").Value = 7991.429
$ws.Range("M77").Value = -3623.429
# Row 102
$ws.Range("H102").Value = 13944.333
$ws.Range("I102").Value = 4425
$ws.Range("K102").Value = 4425
$ws.Range("M102").Value = -2803
# Row 136
$ws.Range("H136").Value = 3065.3823
$ws.Range("I136").Value = 2912.0667
$ws.Range("J136").Value = 4215.25
$ws.Range("K136").Value = 8736.2001
$ws.Range("L136").Value = 12645.75
$ws.Range("M136").Value = -6186.2001
$ws.Range("N136").Value = -17745.75

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 54
$ws.Range("H54").Value = 11331.857
$ws.Range("I54").Value = 11331.857
$ws.Range("K54").Value = 11331.857
$ws.Range("M54").Value = -10847.857
# Row 105
$ws.Range("H105").Value = 1587.591
$ws.Range("I105").Value = 1421.3125
$ws.Range("K105").Value = 1421.3125
$ws.Range("M105").Value = 325.6875
# Row 134
$ws.Range("H134").Value = 1251.0952
$ws.Range("I134").Value = 1109.2222
$ws.Range("J134").Value = 2102.3333
$ws.Range("K134").Value = 3327.6666
$ws.Range("L134").Value = 6306.999899999999
$ws.Range("M134").Value = -792.6665999999996
$ws.Range("N134").Value = -11376.9999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 5245.64
$ws.Range("I31").Value = 3521.7778
$ws.Range("J31").Value = 5624.049
$ws.Range("K31").Value = 3521.7778
$ws.Range("L31").Value = 5624.049
$ws.Range("M31").Value = -3226.7778
$ws.Range("N31").Value = -6214.049
# Row 34
$ws.Range("H34").Value = 5245.64
$ws.Range("I34").Value = 3521.7778
$ws.Range("J34").Value = 5624.049
$ws.Range("K34").Value = 3521.7778
$ws.Range("L34").Value = 5624.049
$ws.Range("M34").Value = -3319.7778
$ws.Range("N34").Value = -6028.049
# Row 69
$ws.Range("H69").Value = 18801
$ws.Range("I69").Value = 5402.5
$ws.Range("K69").Value = 5402.5
$ws.Range("M69").Value = -4653.5
# Row 72
$ws.Range("H72").Value = 18801
$ws.Range("I72").Value = 5402.5
$ws.Range("K72").Value = 16207.5
$ws.Range("M72").Value = -12463.5
# Row 94
$ws.Range("H94").Value = 2921.25
$ws.Range("I94").Value = 2384
$ws.Range("K94").Value = 2384
$ws.Range("M94").Value = -1933
# Row 108
$ws.Range("H108").Value = 59014.75
$ws.Range("J108").Value = 59014.75
$ws.Range("L108").Value = 59014.75
$ws.Range("N108").Value = -66694.75

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 56
$ws.Range("H56").Value = 5920.7144
$ws.Range("I56").Value = 5920.7144
$ws.Range("K56").Value = 5920.7144
$ws.Range("M56").Value = -5390.7144
# Row 60
$ws.Range("H60").Value = 657.75
$ws.Range("I60").Value = 900
$ws.Range("J60").Value = 512.4
$ws.Range("K60").Value = 2700
$ws.Range("L60").Value = 1537.2
$ws.Range("M60").Value = -2449
$ws.Range("N60").Value = -2039.2
# Row 100
$ws.Range("H100").Value = 250
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 250
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 750
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -2372
# Row 118
$ws.Range("H118").Value = 1929.6666
$ws.Range("I118").Value = 894.5
$ws.Range("K118").Value = 2683.5
$ws.Range("M118").Value = -1440.5
# Row 131
$ws.Range("H131").Value = 33827.78
$ws.Range("I131").Value = 168416.67
$ws.Range("J131").Value = 2768.8076
$ws.Range("K131").Value = 505250.01
$ws.Range("L131").Value = 8306.4228
$ws.Range("M131").Value = -500210.01
$ws.Range("N131").Value = -18386.4228

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 22
$ws.Range("H22").Value = 828.3333
$ws.Range("J22").Value = 892.5
$ws.Range("L22").Value = 892.5
$ws.Range("N22").Value = -1950.5
# Row 25
$ws.Range("H25").Value = 895
$ws.Range("J25").Value = 895
$ws.Range("L25").Value = 895
$ws.Range("N25").Value = -1953
# Row 113
$ws.Range("H113").Value = 3706705.8
$ws.Range("I113").Value = 2998.5
$ws.Range("J113").Value = 4764908
$ws.Range("K113").Value = 2998.5
$ws.Range("L113").Value = 4764908
$ws.Range("M113").Value = -828.5
$ws.Range("N113").Value = -4769248
# Row 132
$ws.Range("H132").Value = 2116.394
$ws.Range("I132").Value = 1753.3478
$ws.Range("J132").Value = 2951.4
$ws.Range("K132").Value = 5260.0434
$ws.Range("L132").Value = 8854.200000000001
$ws.Range("M132").Value = -2730.0434
$ws.Range("N132").Value = -13914.2

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 116
$ws.Range("H116").Value = 210999.75
$ws.Range("J116").Value = 210999.75
$ws.Range("L116").Value = 210999.75
$ws.Range("N116").Value = -220177.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 110
$ws.Range("H110").Value = 60644
$ws.Range("J110").Value = 60644
$ws.Range("L110").Value = 60644
$ws.Range("N110").Value = -68824
# Row 133
$ws.Range("H133").Value = 80500
$ws.Range("J133").Value = 80500
$ws.Range("L133").Value = 80500
$ws.Range("N133").Value = -90620
